# aggiornamento fino a 28 luglio
# Append 27 new daily rows (302..328) after the existing last row (301),
# continuing the date series (serial 44375 -> 44402) with zeroed columns
# B/C/D, matching the style used by the existing date column (A).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Style template cell: the last existing date cell (A301), which carries
# the bold/bordered/centered date-number-format style used throughout
# column A.
$styleSource = $ws.Cells.Item(301, 1)

$startRow = 302
$endRow = 328
$startSerial = 44376

for ($i = 0; $i -le ($endRow - $startRow); $i++) {
    $row = $startRow + $i
    $serial = $startSerial + $i

    $dateCell = $ws.Cells.Item($row, 1)
    # Copy formatting (style) from the prior date cell, then set the value.
    $styleSource.Copy($dateCell)
    $dateCell.Value = $serial

    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = 0
}
